$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "G"  = -0.05533980582524266
    "H"  = -0.6766990291262135
    "I"  = -0.6019417475728155
    "J"  = -0.6019417475728155
    "K"  = -6.87
    "L"  = -0.6669902912621359
    "U"  = 6.17
    "V"  = 0.02877798507462687
    "W"  = -0.5325581395348837
    "X"  = 0.05751814610818003
    "Y"  = -0.5900762856430638
    "Z"  = 4.583889630618603
    "AA" = -2.75923453493547
    "AB" = 0.05644456583505471
    "AC" = -2.815679100770525
    "AD" = 5.84
    "AE" = 0
    "AF" = 5.84
    "AG" = -0.3300000000000001
    "AH" = 0.02651652742462768
    "AI" = 0.4329132690882135
    "AJ" = -0.001541551828841034
    "AK" = -0.04508196721311476
    "AL" = 0.345
    "AM" = 0.263
    "AN" = -0.9684908789386401
    "AO" = -17.97101449275362
    "AP" = 0.05472636815920399
    "AQ" = -23.57414448669202
}

foreach ($col in $values.Keys) {
    $newVal = $values[$col]
    $ws.Range("$col" + "2").Value = $newVal
    $ws.Range("$col" + "3").Value = $newVal
}
